# ImageMap.xlsx layout update: credit | jackpot | bet label rows added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (creditImage) / Row 13 (betImage): fill in values that were
#     previously blank placeholder rows ---
$ws.Range("A12").Value = "creditImage"
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = 35
$ws.Range("D12").Value = 180
$ws.Range("E12").Value = 293

$ws.Range("A13").Value = "betImage"
$ws.Range("B13").Value = 80
$ws.Range("C13").Value = 35
$ws.Range("D13").Value = 335
$ws.Range("E13").Value = 293

# --- New rows 14-16: jackpot / credit / bet label ---
$ws.Range("A14").Value = "jackpot label"
$ws.Range("D14").Value = "center X"
$ws.Range("E14").Value = 11

$ws.Range("A15").Value = "credit label"
$ws.Range("D15").Value = 290
$ws.Range("E15").Value = 303

$ws.Range("A16").Value = "bet label"
$ws.Range("D16").Value = 433
$ws.Range("E16").Value = 303

Write-Host "content done"
